$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose values changed (rolling weekly dataset) ---
$ws.Range("D10").Value = 44434
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 24000
$ws.Range("S10").Value = 1200

$ws.Range("D11").Value = 44221
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 25000
$ws.Range("S11").Value = 1250

$ws.Range("D12").Value = 44428
$ws.Range("M12").Value = 15
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 24000
$ws.Range("S12").Value = 1200

$ws.Range("D13").Value = 44175
$ws.Range("M13").Value = 25
$ws.Range("N13").Value = 23000
$ws.Range("O13").Value = 23000
$ws.Range("P13").Value = 23000
$ws.Range("S13").Value = 1150

$ws.Range("D14").Value = 44363
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 24000
$ws.Range("O14").Value = 24000
$ws.Range("P14").Value = 24000
$ws.Range("S14").Value = 1200

$ws.Range("D15").Value = 44349
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 24000
$ws.Range("O15").Value = 24000
$ws.Range("P15").Value = 24000
$ws.Range("S15").Value = 1200

$ws.Range("D16").Value = 44421
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = 24000
$ws.Range("O16").Value = 24000
$ws.Range("P16").Value = 24000
$ws.Range("S16").Value = 1200

$ws.Range("D17").Value = 44222
$ws.Range("M17").Value = 15
$ws.Range("N17").Value = 25000
$ws.Range("O17").Value = 25000
$ws.Range("P17").Value = 25000
$ws.Range("S17").Value = 1250

$ws.Range("D18").Value = 44400
$ws.Range("M18").Value = 5
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 24000
$ws.Range("S18").Value = 1200

$ws.Range("D19").Value = 44426
$ws.Range("M19").Value = 15
$ws.Range("N19").Value = 24000
$ws.Range("O19").Value = 24000
$ws.Range("P19").Value = 24000
$ws.Range("S19").Value = 1200

$ws.Range("D24").Value = 44435
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 24000
$ws.Range("O24").Value = 24000
$ws.Range("P24").Value = 24000
$ws.Range("S24").Value = 1200

$ws.Range("D25").Value = 44431
$ws.Range("M25").Value = 40
$ws.Range("N25").Value = 24000
$ws.Range("O25").Value = 24000
$ws.Range("P25").Value = 24000
$ws.Range("S25").Value = 1200

$ws.Range("D26").Value = 44231
$ws.Range("M26").Value = 15
$ws.Range("N26").Value = 25000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 25000
$ws.Range("S26").Value = 1250

$ws.Range("D27").Value = 44391
$ws.Range("M27").Value = 10
$ws.Range("N27").Value = 24000
$ws.Range("O27").Value = 24000
$ws.Range("P27").Value = 24000
$ws.Range("S27").Value = 1200

$ws.Range("D28").Value = 44389
$ws.Range("M28").Value = 20
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("S28").Value = 1200

$ws.Range("D29").Value = 44251
$ws.Range("M29").Value = 15
$ws.Range("N29").Value = 25000
$ws.Range("O29").Value = 25000
$ws.Range("P29").Value = 25000
$ws.Range("S29").Value = 1250

$ws.Range("D30").Value = 44419
$ws.Range("M30").Value = 40
$ws.Range("N30").Value = 25000
$ws.Range("O30").Value = 25000
$ws.Range("P30").Value = 25000
$ws.Range("S30").Value = 1250

$ws.Range("D31").Value = 44420
$ws.Range("M31").Value = 35
$ws.Range("N31").Value = 25000
$ws.Range("O31").Value = 25000
$ws.Range("P31").Value = 25000
$ws.Range("S31").Value = 1250

$ws.Range("D32").Value = 44433
$ws.Range("M32").Value = 10
$ws.Range("N32").Value = 24000
$ws.Range("O32").Value = 24000
$ws.Range("P32").Value = 24000
$ws.Range("S32").Value = 1200

$ws.Range("D33").Value = 44232
$ws.Range("M33").Value = 15
$ws.Range("N33").Value = 25000
$ws.Range("O33").Value = 25000
$ws.Range("P33").Value = 25000
$ws.Range("S33").Value = 1250

$ws.Range("D34").Value = 44398
$ws.Range("M34").Value = 15
$ws.Range("N34").Value = 25000
$ws.Range("O34").Value = 25000
$ws.Range("P34").Value = 25000
$ws.Range("S34").Value = 1250

# --- Append brand-new rows 35-37 (same constant columns as the rest of the dataset) ---
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44334
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100108
$ws.Range("H35").Value = "Tropicales y subtropicales"
$ws.Range("I35").Value = 100108007
$ws.Range("J35").Value = "Coco"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 20
$ws.Range("N35").Value = 25000
$ws.Range("O35").Value = 25000
$ws.Range("P35").Value = 25000
$ws.Range("Q35").Value = "$/malla 20 unidades"
$ws.Range("R35").Value = "Perú"
$ws.Range("S35").Value = 1250
$ws.Range("T35").Value = 20

$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44418
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108007
$ws.Range("J36").Value = "Coco"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 20
$ws.Range("N36").Value = 24000
$ws.Range("O36").Value = 24000
$ws.Range("P36").Value = 24000
$ws.Range("Q36").Value = "$/malla 20 unidades"
$ws.Range("R36").Value = "Perú"
$ws.Range("S36").Value = 1200
$ws.Range("T36").Value = 20

$ws.Range("A37").Value = 10
$ws.Range("B37").Value = "Vega Modelo de Temuco"
$ws.Range("C37").Value = "La Araucanía"
$ws.Range("D37").Value = 44432
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100108
$ws.Range("H37").Value = "Tropicales y subtropicales"
$ws.Range("I37").Value = 100108007
$ws.Range("J37").Value = "Coco"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 30
$ws.Range("N37").Value = 24000
$ws.Range("O37").Value = 24000
$ws.Range("P37").Value = 24000
$ws.Range("Q37").Value = "$/malla 20 unidades"
$ws.Range("R37").Value = "Perú"
$ws.Range("S37").Value = 1200
$ws.Range("T37").Value = 20
